# Apply cryptos list update (values sourced from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.402.02"
$ws.Range("D3").Value = "1.671.18"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "220.67"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "0.5356"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").Value = "1.011"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "0.06398"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "21.02"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").Value = "0.07861"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "4.563"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.673.73"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "1.900.18"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "0.5562"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "0.0₅8185"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "66.25"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "26.422.12"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'4.690"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "196.42"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "1.012"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "146.24"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "0.1229"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "7.253"
$ws.Range("D28").Value = "16.23"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  +3.70%  "
$ws.Range("D30").Value = "'0.05890"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "3.582"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "1.621"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "0.9713"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("D36").Value = "2.853"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "0.5821"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "1.076.14"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("D41").Value = "0.8677"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").Value = "5.887"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("D44").Value = "104.25"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "1.809.86"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("D46").Value = "58.19"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.011"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈107"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").Value = "0.4401"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").Value = "8.073"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("E51").Value = "  +0.50%  "
